# Update Name of Algo
# Applies updated imputed values to the RandomForest result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = -12.3491
$ws.Range("A9").Value = -20.35539999999998
$ws.Range("C11").Value = -13.673
$ws.Range("A18").Value = -22.90210000000001
$ws.Range("A20").Value = -22.19100000000002
